# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# This script updates the DAMSLTag (column I) and DialogAct (column J) values
# for the rows whose dialog-act annotations changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 3; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 20; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 43; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 44; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 48; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 54; DAMSLTag = '%'; DialogAct = 'Uninterpretable' },
    @{ Row = 57; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 58; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 63; DAMSLTag = '%'; DialogAct = 'Uninterpretable' },
    @{ Row = 68; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 76; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 80; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 84; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 94; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 103; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 104; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 127; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 128; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 134; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 136; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 141; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 150; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 151; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 156; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 163; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 192; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 211; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 215; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 226; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 227; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 228; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 229; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 232; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 237; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 245; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 247; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 248; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 258; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 262; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 263; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 271; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 319; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 339; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 341; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 352; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 359; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 366; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 370; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 374; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 383; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 384; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 387; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 388; DAMSLTag = '%'; DialogAct = 'Uninterpretable' },
    @{ Row = 389; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 390; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 393; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 394; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 395; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 398; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 402; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 403; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 420; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 425; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 434; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 440; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 458; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 464; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 482; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 485; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 503; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 508; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 519; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 529; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 542; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 544; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 546; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 549; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 567; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 572; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 577; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 585; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 590; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 604; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 632; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 637; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 643; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 649; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 650; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 681; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
)

# Column indexes: I = 9 (DAMSLTag), J = 10 (DialogAct)
$damslCol = 9
$dialogActCol = 10

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, $damslCol).Value = $change.DAMSLTag
    $ws.Cells.Item($change.Row, $dialogActCol).Value = $change.DialogAct
}
